# Updated cryptos list (price/volume refresh) on Thu Dec 21 11:54:38 UTC 2023.
# Price/volume cells are stored as plain text in this sheet (not numbers),
# so for any new Price value that looks numeric we briefly force a text
# number format before assigning it (otherwise Excel auto-converts it to a
# real number and strips things like trailing zeros), then restore the
# cell's style so no visible formatting changes remain.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.961.19'
$ws.Range("E2").Value = '  +2.37%  '

$ws.Range("D3").Value = '2.250.67'
$ws.Range("E3").Value = '  +1.60%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '271.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.84%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.613'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.67%  '

$ws.Range("E11").Value = '  +1.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.04%  '

$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("D14").Value = '2.589.32'
$ws.Range("E14").Value = '  +1.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.07%  '

$ws.Range("D16").Value = '2.243.77'
$ws.Range("E16").Value = '  +1.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.801'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.95%  '

$ws.Range("D18").Value = '43.912.67'
$ws.Range("E18").Value = '  +2.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000105'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '

$ws.Range("E20").Value = '  +1.35%  '

$ws.Range("E21").Value = '  -1.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.72%  '

$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.26%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +16.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.20%  '

$ws.Range("E28").Value = '  +6.67%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.92%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '40.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0914'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.20%  '

$ws.Range("E33").Value = '  +2.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.68%  '

$ws.Range("E35").Value = '  +1.99%  '

$ws.Range("E36").Value = '  +5.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0355'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.54%  '

$ws.Range("E43").Value = '  +1.57%  '

$ws.Range("E44").Value = '  +2.15%  '

$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.84%  '

$ws.Range("E48").Value = '  +7.53%  '

$ws.Range("E49").Value = '  +2.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.439'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.12%  '

$ws.Range("E51").Value = '  +1.04%  '
